$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")

# Overview sheet: row 3 corresponds to f2fce4a5-57cd-45c5-ac6d-71a73fa406dc.md
# Status changes from "In Translation" to "Ready for handoff" for both zh-cn (B3) and de-de (C3)
$overview.Range("B3").Value = "Ready for handoff"
$overview.Range("C3").Value = "Ready for handoff"

# zh-cn sheet: row 3 is the f2fce4a5...md row
# Status (B3): "In Translation" -> "Ready for handoff"
# Latest Handoff Datetime (D3): "2016-03-09 02:15:25" -> "2016-03-09 02:17:46"
$zhcn.Range("B3").Value = "Ready for handoff"
$zhcn.Range("D3").Value = "2016-03-09 02:17:46"

# de-de sheet: row 3 is the f2fce4a5...md row
# Status (B3): "In Translation" -> "Ready for handoff"
# Latest Handoff Datetime (D3): "2016-03-09 02:16:15" -> "2016-03-09 02:17:56"
$dede.Range("B3").Value = "Ready for handoff"
$dede.Range("D3").Value = "2016-03-09 02:17:56"
